# "Generate Report for Handoff"
# The localization CI job re-ran and generated a new handoff package for the
# e48d7f06-ea32-4613-9630-985790195de6 source file: its status flips from
# "In Translation" to "Ready for handoff", its priority flips from "ht" to
# "mt", and a fresh handoff timestamp is recorded (per-locale) on both the
# zh-cn and de-de sheets, which the Overview sheet mirrors/rolls up.

$wb = $excel.ActiveWorkbook

$ovw  = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# ---- Overview sheet (row 3 = e48d7f06-...md) ----
$ovw.Range("E3").Value = "Ready for handoff"
$ovw.Range("F3").Value = "Ready for handoff"
$ovw.Range("G3").Value = "2016-08-22 02:11:47"

# ---- zh-cn sheet (row 3 = e48d7f06-...md) ----
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-22 02:11:43"

# ---- de-de sheet (row 3 = e48d7f06-...md) ----
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-22 02:11:47"

# ---- Column widths grow to fit the new, longer "Ready for handoff" text ----
# (Status column on Overview is split across two locale columns, E & F;
#  Status column on the per-locale sheets is column C.)
$ovw.Columns.Item(5).ColumnWidth = 16.333333333333332
$ovw.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
